$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item(1)
$wsMonthly = $wb.Worksheets.Item(2)

# Rename header label on "Weekly Quantity" sheet
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# Rename header label on "Monthly Trend" sheet
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add new "PO Forecast" worksheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "PO Forecast"

# Match page margins used by the other sheets (inches -> points: 0.75"=54pt, 1"=72pt, 0.5"=36pt)
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Headers
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Copy header formatting (bold, centered, bordered) from an existing header cell
$wsWeekly.Range("B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Copy date-style formatting from an existing date cell onto column A (rows 2-12)
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A12").PasteSpecial(-4122)

# Data rows
$newSheet.Range("A2").Value = 45011.99999999999
$newSheet.Range("B2").Value = 10
$newSheet.Range("C2").Value = 9.999999986633208
$newSheet.Range("D2").Value = 10.00000001302823

$newSheet.Range("A3").Value = 45018.99999999999
$newSheet.Range("B3").Value = 10
$newSheet.Range("C3").Value = 9.999999987780257
$newSheet.Range("D3").Value = 10.00000001239442

$newSheet.Range("A4").Value = 45039.99999999999
$newSheet.Range("B4").Value = 10
$newSheet.Range("C4").Value = 9.999999987274562
$newSheet.Range("D4").Value = 10.00000001242842

$newSheet.Range("A5").Value = 45046.99999999999
$newSheet.Range("B5").Value = 10
$newSheet.Range("C5").Value = 9.999999985379104
$newSheet.Range("D5").Value = 10.00000001461711

$newSheet.Range("A6").Value = 45053.99999999999
$newSheet.Range("B6").Value = 10
$newSheet.Range("C6").Value = 9.999999977714417
$newSheet.Range("D6").Value = 10.00000002257774

$newSheet.Range("A7").Value = 45060.99999999999
$newSheet.Range("B7").Value = 10
$newSheet.Range("C7").Value = 9.999999963554613
$newSheet.Range("D7").Value = 10.00000004311949

$newSheet.Range("A8").Value = 45067.99999999999
$newSheet.Range("B8").Value = 10
$newSheet.Range("C8").Value = 9.99999993604415
$newSheet.Range("D8").Value = 10.00000006951655

$newSheet.Range("A9").Value = 45074.99999999999
$newSheet.Range("B9").Value = 10
$newSheet.Range("C9").Value = 9.999999902698981
$newSheet.Range("D9").Value = 10.00000010083982

$newSheet.Range("A10").Value = 45081.99999999999
$newSheet.Range("B10").Value = 10
$newSheet.Range("C10").Value = 9.999999866695717
$newSheet.Range("D10").Value = 10.00000014670086

$newSheet.Range("A11").Value = 45088.99999999999
$newSheet.Range("B11").Value = 10
$newSheet.Range("C11").Value = 9.999999819920445
$newSheet.Range("D11").Value = 10.00000018413291

$newSheet.Range("A12").Value = 45095.99999999999
$newSheet.Range("B12").Value = 10
$newSheet.Range("C12").Value = 9.999999772915922
$newSheet.Range("D12").Value = 10.0000002282208

# Activate the first sheet to preserve the original active tab
$wsWeekly.Activate()
